# Adds missing THT components (J1-J5, J7-J9, SW1, SW2) to the CPL sheet,
# keeping the component list sorted by designator - mirrors a refreshed
# query-table pull from the updated RP2040-Eins-top-pos.csv source.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert J1..J5,J7..J9 (8 new rows) right before the existing J10 row (row 29) ---
$ws.Range("29:36").Insert()

$jrows = @(
    @("J1", 137.572,             -119.38,              0),
    @("J2", 190.495,             -96.52,                0),
    @("J3", 177.8,               -124.46,              90),
    @("J4", 154.95500000000001,  -124.435,             90),
    @("J5", 190.5,               -85.07,              180),
    @("J7", 191.14500000000001,  -106.67,               0),
    @("J8", 190.5,               -76.2,               -90),
    @("J9", 168.626,             -76.224999999999994, -90)
)

$r = 29
foreach ($row in $jrows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = "top"
    $ws.Range("E$r").Value = $row[3]
    $r = $r + 1
}

# --- Insert SW1, SW2 (2 new rows) right after R15, before U1 (row 55 after the above insert) ---
$ws.Range("55:56").Insert()

$swrows = @(
    @("SW1", 168.91, -107.95,   0),
    @("SW2", 133.35,  -77.47, 180)
)

$r = 55
foreach ($row in $swrows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = "top"
    $ws.Range("E$r").Value = $row[3]
    $r = $r + 1
}

# --- Update the defined name range to cover the new extent ---
$wb.Names.Item("Sheet1!RP2040_Eins_top_pos").RefersTo = "=Sheet1!`$A`$1:`$E`$62"

# --- Update the header selection / scroll position to match the saved view ---
$ws.Range("A1:E1").Select()
